$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the comma in the "Context" column (column I) values for the
# burial/midden rows - e.g. "burial, M020" -> "burial M020".
$ws.Range("I10").Value = "burial M020"
$ws.Range("I11").Value = "midden H044"
$ws.Range("I12").Value = "burial M009"
$ws.Range("I16").Value = "burial M039"
$ws.Range("I19").Value = "burial M066"
$ws.Range("I21").Value = "midden H193"
$ws.Range("I23").Value = "midden H026"
$ws.Range("I24").Value = "burial M095"

# Move the active selection / scroll position to reflect where the author
# was working when the file was last saved.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("K12").Select()
